$d = $word.ActiveDocument
$footer = $d.Sections(1).Footers(1)

# --- Paragraph 2: "Centro Médico: Av. Rumiñahui 3656 y Roosvelt"
#     -> "Consultorio: «=consultation.branch_office.address»"
$p2xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="58219ADA" w14:textId="77777777" w:rsidR="004D2A29" w:rsidRPr="006D7442" w:rsidRDefault="004D2A29" w:rsidP="004D2A29"><w:pPr><w:pStyle w:val="Footer"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Consultorio:</w:t></w:r><w:r w:rsidRPr="006D7442"><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD =consultation.branch_office.address \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>«=consultation.branch_office.address»</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>'
$footer.Range.Paragraphs(2).Range.InsertXML($p2xml)

# --- Paragraph 3: "Cita Previa: (03) 280 0336 – (03) 281 4219 / Emergencias: 098 458 4566"
#     -> "Cita Previa: «=consultation.branch_office.phone_number» / Emergencias: «=emergency_number»"
$p3xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1C3E8114" w14:textId="77777777" w:rsidR="004D2A29" w:rsidRPr="006D7442" w:rsidRDefault="004D2A29" w:rsidP="004D2A29"><w:pPr><w:pStyle w:val="Footer"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="006D7442"><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Cita Previa:</w:t></w:r><w:r w:rsidRPr="006D7442"><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD =consultation.branch_office.phone_numbers \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>«=consultation.branch_office.phone_number»</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> / </w:t></w:r><w:r w:rsidRPr="006D7442"><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Emergencias:</w:t></w:r><w:r w:rsidRPr="006D7442"><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD =emergency_number \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>«=emergency_number»</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>'
$footer.Range.Paragraphs(3).Range.InsertXML($p3xml)

# --- Paragraph 5: "Latacunga-Ecuador" -> "«=consultation.branch_office.city»-Ecuador"
$p5xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5654154A" w14:textId="77777777" w:rsidR="004D2A29" w:rsidRPr="006D7442" w:rsidRDefault="004D2A29" w:rsidP="004D2A29"><w:pPr><w:pStyle w:val="Footer"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="006D7442"><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD =consultation.branch_office.city \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>«=consultation.branch_office.city»</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>-Ecuador</w:t></w:r></w:p>'
$footer.Range.Paragraphs(5).Range.InsertXML($p5xml)

Write-Output "Footer final text: [$($footer.Range.Text)]"
